$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed) date column C for rows 2-7 from 45184 to 45186
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45186
}

# Row 2 (A 55942-2019): add friendly name as second HYPERLINK argument
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_KAVLINGE/artfynd/A 55942-2019.xlsx", "A 55942-2019")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_KAVLINGE/kartor/A 55942-2019.png", "A 55942-2019")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_KAVLINGE/klagomål/A 55942-2019.docx", "A 55942-2019")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_KAVLINGE/klagomålsmail/A 55942-2019.docx", "A 55942-2019")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_KAVLINGE/tillsyn/A 55942-2019.docx", "A 55942-2019")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_KAVLINGE/tillsynsmail/A 55942-2019.docx", "A 55942-2019")'

# Row 3 (A 1481-2022): add friendly name as second HYPERLINK argument
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_KAVLINGE/artfynd/A 1481-2022.xlsx", "A 1481-2022")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_KAVLINGE/kartor/A 1481-2022.png", "A 1481-2022")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_KAVLINGE/klagomål/A 1481-2022.docx", "A 1481-2022")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_KAVLINGE/klagomålsmail/A 1481-2022.docx", "A 1481-2022")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_KAVLINGE/tillsyn/A 1481-2022.docx", "A 1481-2022")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_KAVLINGE/tillsynsmail/A 1481-2022.docx", "A 1481-2022")'
